$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: Read our review of Amazon Queen..."
#    paragraph that currently sits right under the title heading.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
    }
}

# ------------------------------------------------------------------
# 2. At the end of the document, the paragraph that used to hold the
#    feature-image prompt ("Create a feature image for Amazon
#    Queen: ...") is replaced by two paragraphs:
#      - a bold paragraph repeating the page title
#      - an italic paragraph with the (moved) meta-description text
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$replacementXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Amazon Queen for Free - Slot Game Review</w:t></w:r></w:p>" +
                   "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Amazon Queen, a high-volatility slot with modest betting range, and play for free.</w:t></w:r></w:p>"

$lastPara.Range.InsertXML($replacementXml) | Out-Null
